$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "program" rule body: (expr;)* -> (stmt)*
$ws.Range("C4").Value = "(stmt)*"

# Rename the array-index grammar rule: expr[epxr] -> expr[epxr_list]
$ws.Range("C29").Value = "expr[epxr_list]"

# New example cell next to the if/else rule row: a sample while-condition
$ws.Range("F32").Value = "while (bool1 && bool2)"

# Move the selection to the new cell (matches the saved cursor position)
$ws.Range("F32").Select()
